# FUNCTIONALITY: Re-wrote the Suite around two new test cases.
#
# The underlying test-suite statistics sheet tracks, per item (row), an
# "Automated Test Cases" count (col B) and a "Total Test Cases" count
# (col C). Two additional test cases were written for the row identified
# as "ImportSettingsCardUI" (row 5): both its automated-case count and its
# total-case count go up by 2 (4 -> 6). The roll-up formulas in column G
# (which SUM down columns B and C, and the percentage-automated ratio)
# recalculate automatically from that single edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5 = "ImportSettingsCardUI": two new automated test cases were added.
$ws.Range("B5").Value = 6
$ws.Range("C5").Value = 6
